# consensus flowchart.pptx — wording fix
#
# "Received other agent's solution?" -> "Received other agents' solution?"
#
# The diamond/decision shape ("Losango 11") on slide 1 asks whether the
# other agent's solution has been received. The apostrophe moves from a
# singular possessive ("agent's") to a plural possessive ("agents'").
#
# We edit via TextRange.Characters(start,len) sub-ranges (rather than
# replacing the whole TextRange.Text) so that only the two runs that
# actually changed are touched and the rest of the paragraph's run
# layout / formatting (bold "b=1", err="1" spell-flags, etc.) is left
# exactly as PowerPoint would leave it for a narrow, in-place edit.

$p = $ppt.ActivePresentation

$target = $null
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    $hasText = $false
    try { $hasText = $shp.HasTextFrame } catch { $hasText = $false }
    if (-not $hasText) {
        try { $hasText = ($shp.TextFrame.HasText -ne 0) } catch { $hasText = $false }
    }
    if ($hasText) {
        $t = ""
        try { $t = $shp.TextFrame.TextRange.Text } catch { $t = "" }
        if ($t -like "*agent's solution?*") {
            $target = $shp
            break
        }
    }
}

if ($target -ne $null) {
    # The diamond shape has <a:spAutoFit/>, so rewriting any text inside
    # it makes the host recompute the shape's box. The real edit never
    # touched the shape's geometry (the diff has no <a:xfrm> change), so
    # snapshot the exact position/size first and restore it afterwards
    # (full float precision to avoid EMU round-trip drift).
    $origLeft   = $target.Left
    $origTop    = $target.Top
    $origWidth  = $target.Width
    $origHeight = $target.Height

    $tr = $target.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf("agent's")

    if ($idx -ge 0) {
        # 1-based character position of the "agent's" run.
        $start = $idx + 1

        # "agent's" (7 chars) -> "agents" (drop the apostrophe)
        $word = $tr.Characters($start, 7)
        $word.Text = "agents"

        # the single space right after now sits right after "agents";
        # turn it into "' " (apostrophe + space)
        $space = $tr.Characters($start + 6, 1)
        $space.Text = [string]([char]0x2019) + " "
    }

    $target.Left   = $origLeft
    $target.Top    = $origTop
    $target.Width  = $origWidth
    $target.Height = $origHeight
}

Write-Output "done"
